$d = $word.ActiveDocument

# 1. {GESTOR_CELULAR}} -> {GESTOR_CELULAR}
$d.Content.Find.Execute("{GESTOR_CELULAR}}", $true, $false, $false, $false, $false, $true, 1, $false, "{GESTOR_CELULAR}", 2)

# 2. {GERENTE_CELULAR}} -> {GERENTE_CELULAR}
$d.Content.Find.Execute("{GERENTE_CELULAR}}", $true, $false, $false, $false, $false, $true, 1, $false, "{GERENTE_CELULAR}", 2)

# 3. merge the two runs "{" + "GERENTE}" into a single run "{GERENTE}"
$d.Content.Find.Execute("{GERENTE}", $true, $false, $false, $false, $false, $true, 1, $false, "{GERENTE}", 2)
